$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-data refresh to the per-job
# "Profits" sheets: recomputed average market prices (H/I/J), leve
# totals (K/L) and resulting profit figures (M/N) for the affected rows.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 3710.6667
$ws.Range("I69").Value = 3500
$ws.Range("J69").Value = 3743.077
$ws.Range("K69").Value = 10500
$ws.Range("L69").Value = 11229.231
$ws.Range("M69").Value = -9626
$ws.Range("N69").Value = -12977.231
# Row 72
$ws.Range("H72").Value = 3710.6667
$ws.Range("I72").Value = 3500
$ws.Range("J72").Value = 3743.077
$ws.Range("K72").Value = 31500
$ws.Range("L72").Value = 33687.693
$ws.Range("M72").Value = -27132
$ws.Range("N72").Value = -42423.693
# Row 137
$ws.Range("H137").Value = 27330.395
$ws.Range("I137").Value = 678.9583
$ws.Range("J137").Value = 73018.57000000001
$ws.Range("K137").Value = 2036.8749
$ws.Range("L137").Value = 219055.71
$ws.Range("M137").Value = 513.1251
$ws.Range("N137").Value = -224155.71

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6786.18
$ws.Range("I32").Value = 5104.2046
$ws.Range("J32").Value = 14998.177
$ws.Range("K32").Value = 5104.2046
$ws.Range("L32").Value = 14998.177
$ws.Range("M32").Value = -4817.2046
$ws.Range("N32").Value = -15572.177
# Row 56
$ws.Range("H56").Value = 16217.875
$ws.Range("J56").Value = 16217.875
$ws.Range("L56").Value = 16217.875
$ws.Range("N56").Value = -17701.875
# Row 61
$ws.Range("H61").Value = 1934.3541
$ws.Range("I61").Value = 1672.85
$ws.Range("J61").Value = 3241.875
$ws.Range("K61").Value = 1672.85
$ws.Range("L61").Value = 3241.875
$ws.Range("M61").Value = -1460.85
$ws.Range("N61").Value = -3665.875
# Row 74
$ws.Range("H74").Value = 559555.8
$ws.Range("I74").Value = 5265.9165
$ws.Range("J74").Value = 1668135.6
$ws.Range("K74").Value = 5265.9165
$ws.Range("L74").Value = 1668135.6
$ws.Range("M74").Value = -4391.9165
$ws.Range("N74").Value = -1669883.6
# Row 77
$ws.Range("H77").Value = 559555.8
$ws.Range("I77").Value = 5265.9165
$ws.Range("J77").Value = 1668135.6
$ws.Range("K77").Value = 26329.5825
$ws.Range("L77").Value = 8340678
$ws.Range("M77").Value = -21961.5825
$ws.Range("N77").Value = -8349414
# Row 88
$ws.Range("H88").Value = 166801420
$ws.Range("J88").Value = 200161300
$ws.Range("L88").Value = 200161300
$ws.Range("N88").Value = -200162112
# Row 91
$ws.Range("H91").Value = 166801420
$ws.Range("J91").Value = 200161300
$ws.Range("L91").Value = 200161300
$ws.Range("N91").Value = -200164108
# Row 102
$ws.Range("H102").Value = 1199
$ws.Range("I102").Value = 1148.625
$ws.Range("K102").Value = 1148.625
$ws.Range("M102").Value = 473.375
# Row 110
$ws.Range("H110").Value = 1030.1482
$ws.Range("I110").Value = 900.63635
$ws.Range("J110").Value = 1600
$ws.Range("K110").Value = 900.63635
$ws.Range("L110").Value = 1600
$ws.Range("M110").Value = 1144.36365
$ws.Range("N110").Value = -5690
# Row 122
$ws.Range("H122").Value = 1976.1455
$ws.Range("I122").Value = 1912.1714
$ws.Range("J122").Value = 2088.1
$ws.Range("K122").Value = 5736.5142
$ws.Range("L122").Value = 6264.299999999999
$ws.Range("M122").Value = -3286.5142
$ws.Range("N122").Value = -11164.3
# Row 132
$ws.Range("H132").Value = 13363.917
$ws.Range("I132").Value = 16193.424
$ws.Range("J132").Value = 2989.0557
$ws.Range("K132").Value = 48580.272
$ws.Range("L132").Value = 8967.167099999999
$ws.Range("M132").Value = -46050.272
$ws.Range("N132").Value = -14027.1671
# Row 136
$ws.Range("H136").Value = 1934.3541
$ws.Range("I136").Value = 1672.85
$ws.Range("J136").Value = 3241.875
$ws.Range("K136").Value = 5018.549999999999
$ws.Range("L136").Value = 9725.625
$ws.Range("M136").Value = -2468.549999999999
$ws.Range("N136").Value = -14825.625

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1972.381
$ws.Range("I105").Value = 1624.6154
$ws.Range("J105").Value = 2537.5
$ws.Range("K105").Value = 1624.6154
$ws.Range("L105").Value = 2537.5
$ws.Range("M105").Value = 122.3846000000001
$ws.Range("N105").Value = -6031.5
# Row 107
$ws.Range("H107").Value = 3393
$ws.Range("I107").Value = 4172.7666
$ws.Range("J107").Value = 1266.3636
$ws.Range("K107").Value = 4172.7666
$ws.Range("L107").Value = 1266.3636
$ws.Range("M107").Value = -2252.7666
$ws.Range("N107").Value = -5106.3636
# Row 134
$ws.Range("H134").Value = 223355
$ws.Range("I134").Value = 294641.1
$ws.Range("J134").Value = 3016.182
$ws.Range("K134").Value = 883923.2999999999
$ws.Range("L134").Value = 9048.545999999998
$ws.Range("M134").Value = -881388.2999999999
$ws.Range("N134").Value = -14118.546

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1383.7407
$ws.Range("I31").Value = 932.1613
$ws.Range("J31").Value = 2857.3157
$ws.Range("K31").Value = 932.1613
$ws.Range("L31").Value = 2857.3157
$ws.Range("M31").Value = -637.1613
$ws.Range("N31").Value = -3447.3157
# Row 34
$ws.Range("H34").Value = 1383.7407
$ws.Range("I34").Value = 932.1613
$ws.Range("J34").Value = 2857.3157
$ws.Range("K34").Value = 932.1613
$ws.Range("L34").Value = 2857.3157
$ws.Range("M34").Value = -730.1613
$ws.Range("N34").Value = -3261.3157
# Row 58
$ws.Range("H58").Value = 1851.6428
$ws.Range("I58").Value = 1263.5264
$ws.Range("J58").Value = 3093.2222
$ws.Range("K58").Value = 1263.5264
$ws.Range("L58").Value = 3093.2222
$ws.Range("M58").Value = -1060.5264
$ws.Range("N58").Value = -3499.2222
# Row 62
$ws.Range("H62").Value = 3008.5881
$ws.Range("I62").Value = 2466.6667
$ws.Range("J62").Value = 3304.182
$ws.Range("K62").Value = 2466.6667
$ws.Range("L62").Value = 3304.182
$ws.Range("M62").Value = -1842.6667
$ws.Range("N62").Value = -4552.182
# Row 65
$ws.Range("H65").Value = 3008.5881
$ws.Range("I65").Value = 2466.6667
$ws.Range("J65").Value = 3304.182
$ws.Range("K65").Value = 12333.3335
$ws.Range("L65").Value = 16520.91
$ws.Range("M65").Value = -9213.333500000001
$ws.Range("N65").Value = -22760.91
# Row 94
$ws.Range("H94").Value = 1333.1666
$ws.Range("J94").Value = 999.8
$ws.Range("L94").Value = 999.8
$ws.Range("N94").Value = -1901.8
# Row 134
$ws.Range("H134").Value = 2777.3103
$ws.Range("I134").Value = 2749.8223
$ws.Range("J134").Value = 2872.4614
$ws.Range("K134").Value = 8249.466899999999
$ws.Range("L134").Value = 8617.3842
$ws.Range("M134").Value = -5714.466899999999
$ws.Range("N134").Value = -13687.3842
# Row 136
$ws.Range("H136").Value = 1851.6428
$ws.Range("I136").Value = 1263.5264
$ws.Range("J136").Value = 3093.2222
$ws.Range("K136").Value = 3790.5792
$ws.Range("L136").Value = 9279.6666
$ws.Range("M136").Value = -1240.5792
$ws.Range("N136").Value = -14379.6666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 5421.684
$ws.Range("I70").Value = 3004
$ws.Range("J70").Value = 5875
$ws.Range("K70").Value = 9012
$ws.Range("L70").Value = 17625
$ws.Range("M70").Value = -8697
$ws.Range("N70").Value = -18255
# Row 73
$ws.Range("H73").Value = 5421.684
$ws.Range("I73").Value = 3004
$ws.Range("J73").Value = 5875
$ws.Range("K73").Value = 9012
$ws.Range("L73").Value = 17625
$ws.Range("M73").Value = -7920
$ws.Range("N73").Value = -19809

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1436.8823
$ws.Range("J97").Value = 1251.7
$ws.Range("L97").Value = 1251.7
$ws.Range("N97").Value = -2243.7
# Row 102
$ws.Range("H102").Value = 879.6667
$ws.Range("I102").Value = 739.25806
$ws.Range("J102").Value = 1423.75
$ws.Range("K102").Value = 739.25806
$ws.Range("L102").Value = 1423.75
$ws.Range("M102").Value = 882.74194
$ws.Range("N102").Value = -4667.75
# Row 113
$ws.Range("H113").Value = 1966.6666
$ws.Range("I113").Value = 1900
$ws.Range("K113").Value = 1900
$ws.Range("M113").Value = 270
# Row 132
$ws.Range("H132").Value = 2061.1296
$ws.Range("I132").Value = 1459.9
$ws.Range("J132").Value = 3778.9285
$ws.Range("K132").Value = 4379.700000000001
$ws.Range("L132").Value = 11336.7855
$ws.Range("M132").Value = -1849.700000000001
$ws.Range("N132").Value = -16396.7855

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1574.2
$ws.Range("I40").Value = 1350.8889
$ws.Range("J40").Value = 1909.1666
$ws.Range("K40").Value = 1350.8889
$ws.Range("L40").Value = 1909.1666
$ws.Range("M40").Value = -1214.8889
$ws.Range("N40").Value = -2181.1666
# Row 87
$ws.Range("H87").Value = 10000
$ws.Range("I87").Value = 10000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 10000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -8877
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 10000
$ws.Range("I90").Value = 10000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 30000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -24384
$ws.Range("N90").ClearContents()
# Row 93
$ws.Range("H93").Value = 1193.8334
$ws.Range("I93").Value = 1209.875
$ws.Range("J93").Value = 1161.75
$ws.Range("K93").Value = 1209.875
$ws.Range("L93").Value = 1161.75
$ws.Range("M93").Value = 38.125
$ws.Range("N93").Value = -3657.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 181276.47
$ws.Range("I62").Value = 5166.6665
$ws.Range("J62").Value = 379400
$ws.Range("K62").Value = 5166.6665
$ws.Range("L62").Value = 379400
$ws.Range("M62").Value = -4542.6665
$ws.Range("N62").Value = -380648
# Row 65
$ws.Range("H65").Value = 181276.47
$ws.Range("I65").Value = 5166.6665
$ws.Range("J65").Value = 379400
$ws.Range("K65").Value = 25833.3325
$ws.Range("L65").Value = 1897000
$ws.Range("M65").Value = -22713.3325
$ws.Range("N65").Value = -1903240
# Row 81
$ws.Range("H81").Value = 40001690
$ws.Range("I81").Value = 58824740
$ws.Range("J81").Value = 2699.75
$ws.Range("K81").Value = 117649480
$ws.Range("L81").Value = 5399.5
$ws.Range("M81").Value = -117648419
$ws.Range("N81").Value = -7521.5
# Row 84
$ws.Range("H84").Value = 40001690
$ws.Range("I84").Value = 58824740
$ws.Range("J84").Value = 2699.75
$ws.Range("K84").Value = 588247400
$ws.Range("L84").Value = 26997.5
$ws.Range("M84").Value = -588242096
$ws.Range("N84").Value = -37605.5
# Row 113
$ws.Range("H113").Value = 83333940
$ws.Range("I113").Value = 300
$ws.Range("K113").Value = 900
$ws.Range("M113").Value = 1270
# Row 132
$ws.Range("H132").Value = 1731.2972
$ws.Range("I132").Value = 889.44446
$ws.Range("J132").Value = 2528.842
$ws.Range("K132").Value = 2668.33338
$ws.Range("L132").Value = 7586.526
$ws.Range("M132").Value = -138.33338
$ws.Range("N132").Value = -12646.526
# Row 136
$ws.Range("H136").Value = 7931.1875
$ws.Range("I136").Value = 871.4286
$ws.Range("J136").Value = 13422.111
$ws.Range("K136").Value = 2614.2858
$ws.Range("L136").Value = 40266.333
$ws.Range("M136").Value = -64.28579999999965
$ws.Range("N136").Value = -45366.333

